$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 20, pushing the existing rows 20-33 down to 22-35.
$ws.Rows("20:21").Insert()

# New row 20: Papaya, "Primera" quality, dated 2021-08-05, $/kilo (en caja de 15 kilos)
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44413
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100108
$ws.Range("H20").Value = "Tropicales y subtropicales"
$ws.Range("I20").Value = 100108004
$ws.Range("J20").Value = "Papaya"
$ws.Range("K20").Value = "Cultivar IV Región"
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = 2500
$ws.Range("O20").Value = 2500
$ws.Range("P20").Value = 2500
$ws.Range("Q20").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R20").Value = "Provincia del Elquí"
$ws.Range("S20").Value = 2500
$ws.Range("T20").Value = 1

# New row 21: Papaya, "Segunda" quality, dated 2021-08-05, $/kilo (en caja de 15 kilos)
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 44413
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100108
$ws.Range("H21").Value = "Tropicales y subtropicales"
$ws.Range("I21").Value = 100108004
$ws.Range("J21").Value = "Papaya"
$ws.Range("K21").Value = "Cultivar IV Región"
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 55
$ws.Range("N21").Value = 2000
$ws.Range("O21").Value = 2000
$ws.Range("P21").Value = 2000
$ws.Range("Q21").Value = "`$/kilo (en caja de 15 kilos)"
$ws.Range("R21").Value = "Provincia del Elquí"
$ws.Range("S21").Value = 2000
$ws.Range("T21").Value = 1

# Make sure the style used for date cells (numeric date format) is applied
# to the D column on the two newly inserted rows (mirrors the rest of column D)
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D21").NumberFormat = "YYYY-MM-DD HH:MM:SS"
